# Generate Report for Handback
# Regenerates the handback-status report with a fresh pair of file GUIDs
# (216dc488-9440-4e8f-9088-079353c20f12 / ffff6b7d147c-24ec-4e6b-995c-6b8813252c28)
# and fresh handoff/handback timestamps, replacing the previous run's values
# (841ee786-d700-45b2-b10e-c5ba1b54794e / ffcd7771-7bc1-4e37-9551-009a5ad17b8e).

$wb = $excel.ActiveWorkbook

$oldGuid1 = "841ee786-d700-45b2-b10e-c5ba1b54794e"
$oldGuid2 = "ffcd7771-7bc1-4e37-9551-009a5ad17b8e"
$newGuid1 = "216dc488-9440-4e8f-9088-079353c20f12"
$newGuid2 = "ffff6b7d147c-24ec-4e6b-995c-6b8813252c28"

$oldHash1zh = "ce2899ff474770c16dbd9c2ffa3ae44391a0294a"
$oldHash2zh = "2eb7e006e768dc1613eaa5329f7186761c67984f"
$newHash = "482de819d50684525c7ea94102c53a270ef3d1b8"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid1.md"
$wsOverview.Range("G2").Value = "2016-08-14 03:32:04"

$wsOverview.Range("A3").Value = "$newGuid2.md"
$wsOverview.Range("G3").Value = "2016-08-14 03:32:04"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/7bba1a0795406b6b8fbc36761e881abd56b58d0f/e2e/$oldGuid1.md", "", "", "e2e\$newGuid1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/7bba1a0795406b6b8fbc36761e881abd56b58d0f/e2e/$oldGuid2.md", "", "", "e2e\$newGuid2.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("G2").Value = "$newGuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-14 03:31:54"
$wsZhCn.Range("J2").Value = "$newGuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-14 03:32:26"

$wsZhCn.Range("G3").Value = "$newGuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-14 03:31:54"
$wsZhCn.Range("J3").Value = "$newGuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-14 03:32:26"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/7bba1a0795406b6b8fbc36761e881abd56b58d0f/e2e/$oldGuid1.md", "", "", "$newGuid1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f808a37b88e5800efcdec4d96624dcacf46f4d5e/e2e/$oldGuid1.md", "", "", "$newGuid1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/7bba1a0795406b6b8fbc36761e881abd56b58d0f/e2e/$oldGuid2.md", "", "", "$newGuid2.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f808a37b88e5800efcdec4d96624dcacf46f4d5e/e2e/$oldGuid2.md", "", "", "$newGuid2.md")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("G2").Value = "$newGuid1.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-14 03:32:04"
$wsDeDe.Range("J2").Value = "$newGuid1.$newHash.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-14 03:32:36"

$wsDeDe.Range("G3").Value = "$newGuid1.$newHash.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-14 03:32:04"
$wsDeDe.Range("J3").Value = "$newGuid1.$newHash.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-14 03:32:36"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/7bba1a0795406b6b8fbc36761e881abd56b58d0f/e2e/$oldGuid1.md", "", "", "$newGuid1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/729dd6a307369b32860d84250f0348db55a51acd/e2e/$oldGuid1.md", "", "", "$newGuid1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/7bba1a0795406b6b8fbc36761e881abd56b58d0f/e2e/$oldGuid2.md", "", "", "$newGuid2.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/729dd6a307369b32860d84250f0348db55a51acd/e2e/$oldGuid2.md", "", "", "$newGuid2.md")
